# Fruta / hortaliza, semanal
# Insert a new weekly record at row 86 (pushing the existing rows 86-122 down
# to 87-123) in the "Vega Modelo de Temuco - Arándano (blue)" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 86..122 down to 87..123 by inserting a blank row at 86.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the new weekly record.
$ws.Range("A86").Value = 10
$ws.Range("B86").Value = "Vega Modelo de Temuco"
$ws.Range("C86").Value = "La Araucanía"
$ws.Range("D86").Value = 44917
$ws.Range("E86").Value = 9
$ws.Range("F86").Value = "Fruta"
$ws.Range("G86").Value = 100101
$ws.Range("H86").Value = "Berries"
$ws.Range("I86").Value = 100101001
$ws.Range("J86").Value = "Arándano (blue)"
$ws.Range("K86").Value = "Sin especificar"
$ws.Range("L86").Value = "Primera"
$ws.Range("M86").Value = 600
$ws.Range("N86").Value = 2000
$ws.Range("O86").Value = 2000
$ws.Range("P86").Value = 2000
$ws.Range("Q86").Value = "$/kilo"
$ws.Range("R86").Value = "Región del Maule"
$ws.Range("S86").Value = 2000
$ws.Range("T86").Value = 1
